# B6-PowerPoint.pptx edit
#
# 1) Three tables (slides 14, 15, 16) get switched from the deck's custom
#    "Table_0" style to the built-in "Medium Style 2 - Accent 1" table
#    style ({C4D6CF68-3FBA-481B-B749-F9C22752F284}).
# 2) The presentation's theme colour palette is reset from the "Integral"
#    (Red Violet) palette back to the default Office palette - i.e. the
#    deck's Design goes back to plain "Office Theme" colours.

$p = $ppt.ActivePresentation

# --- 1) Table styles ---------------------------------------------------
$newTableStyle = "{C4D6CF68-3FBA-481B-B749-F9C22752F284}"
foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newTableStyle)
        }
    }
}

# --- 2) Theme colours: Integral / Red Violet -> Office -----------------
$theme = $p.Slides.Item(1).Master.Theme
$tcs = $theme.ThemeColorScheme

$tcs.Colors(1).RGB  = 0         # dk1      000000
$tcs.Colors(2).RGB  = 16777215  # lt1      FFFFFF
$tcs.Colors(3).RGB  = 6968388   # dk2      44546A
$tcs.Colors(4).RGB  = 15132391  # lt2      E7E6E6
$tcs.Colors(5).RGB  = 13998939  # accent1  5B9BD5
$tcs.Colors(6).RGB  = 3243501   # accent2  ED7D31
$tcs.Colors(7).RGB  = 10855845  # accent3  A5A5A5
$tcs.Colors(8).RGB  = 49407     # accent4  FFC000
$tcs.Colors(9).RGB  = 12874308  # accent5  4472C4
$tcs.Colors(10).RGB = 4697456   # accent6  70AD47
$tcs.Colors(11).RGB = 12673797  # hlink    0563C1
$tcs.Colors(12).RGB = 7491477   # folHlink 954F72
